$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing descriptions with more detail ---
$ws.Range("C3").Value = "Allows data to be stored and retrieved from external serial FLASH memory. There are 8 FLASH chips which each hold a slide of data for the display"
$ws.Range("C4").Value = "Serial connection for WIFI module, used for moving image data from an android app into the embedded system"
$ws.Range("C5").Value = "Serial connection for USB debugging, virtual COM port"

# --- New row 6: Watchdog Timer ---
$ws.Range("A6").Value = "WDT"
$ws.Range("B6").Value = "Watchdog Timer"
$ws.Range("C6").Value = "Error catch, resets the microcontroller if the timer has not been cleared in 2 seconds. Timer is cleared within the heartbeat ISR (below)"
$ws.Range("D6").Value = "watchdogTimer"
$ws.Range("E6").Value = "Postscaler set to 2.048 seconds"

# --- New row 7: Timer 1 (heartbeat timer) ---
$ws.Range("A7").Value = "TMR1"
$ws.Range("B7").Value = "Timer 1"
$ws.Range("C7").Value = "Timer 1 is used as the heartbeat timer, which triggers an interrupt every second. This blinks an LED, clears the watchdog timer, and increments on time counters. Used as a 1 Hz timebase"
$ws.Range("D7").Value = "heartbeatTimer"
$ws.Range("E7").Value = "Timer 1 period set to 61523, input frequency is 15.75 MHz, yields an IRQ rate of 1 Hz"

# Carry the wrapped-text body style from an existing data row onto the new rows
$ws.Range("A2:E2").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122) # xlPasteFormats

# Row heights grow to fit the longer wrapped text (matches the values Excel
# computed for the new/updated rows)
$ws.Rows.Item(3).RowHeight = 72.5
$ws.Rows.Item(4).RowHeight = 58
$ws.Rows.Item(6).RowHeight = 72.5
$ws.Rows.Item(7).RowHeight = 87

# Selection marker mirrors the author's last click position after editing
$ws.Range("E8").Select()
